$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 666, shifting existing rows 666..727 down to 667..728.
$ws.Rows.Item(666).Insert()

# Populate the newly inserted row 666 with the new record's data.
$ws.Range("A666").Value = 5
$ws.Range("B666").Value = "Macroferia Regional de Talca"
$ws.Range("C666").Value = "Maule"
$ws.Range("D666").Value = 45166
$ws.Range("E666").Value = 7
$ws.Range("F666").Value = 100112043
$ws.Range("G666").Value = "Pepino ensalada"
$ws.Range("H666").Value = "Sin especificar"
$ws.Range("I666").Value = "Primera"
$ws.Range("J666").Value = 500
$ws.Range("K666").Value = 7000
$ws.Range("L666").Value = 7000
$ws.Range("M666").Value = 7000
$ws.Range("N666").Value = "`$/caja 60 unidades"
$ws.Range("O666").Value = "Región de Arica y Parinacota"
$ws.Range("P666").Value = 117
$ws.Range("Q666").Value = 60
$ws.Range("R666").Value = "Hortaliza"
